$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 - append new data row beneath the existing header/data rows.
# Column A holds a date-like string ("2025-10-17") that must stay a literal
# text value (matching the existing row 2 pattern) instead of being
# auto-converted into a date serial number. Force text format before the
# assignment, then clear the formatting again so the cell is left with no
# explicit style (matching the source data, which carries no style index).
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2025-10-17"
$ws.Range("A3").ClearFormats()

$ws.Range("B3").Value = "Ycc"
$ws.Range("C3").Value = "123ABX000"
$ws.Range("D3").Value = "Kar"
$ws.Range("E3").Value = "Chennai"
